# step_5 updated and reran
# Adds two new worksheets:
#   - all_drugs_present : union of all rows from the four existing
#     per-group sheets, sorted by the original row index (column A)
#   - not_in_DiSCoVER   : small lookup sheet for drugs not found in DiSCoVER
#
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("G3_effective")
$ws2 = $wb.Worksheets.Item("G3_ineffective")
$ws3 = $wb.Worksheets.Item("notG3_effective")
$ws4 = $wb.Worksheets.Item("SHH_effective")

# ---------------------------------------------------------------
# 1) all_drugs_present
# ---------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add($null, $lastSheet)
$ws5.Name = "all_drugs_present"

# Header row - copy formatting (bold + border + centered) from an existing
# sheet's header row, then (re)write the shared-string values.
$ws1.Range("B1:L1").Copy($ws5.Range("B1:L1"))

$headers = @("Name","G3","mean_score(G3)","G4","mean_score(G4)","SHH","mean_score(SHH)","SHH+p53","mean_score(SHH+p53)","SBI","sbi_id")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws5.Cells.Item(1, 2 + $c).Value = $headers[$c]
}

# Data rows - the union of all rows from the 4 source sheets, sorted by the
# original dataframe index (column A).
$rows = @(
    @{A=1;   Src=$ws1; SrcRow=2},
    @{A=5;   Src=$ws3; SrcRow=2},
    @{A=7;   Src=$ws1; SrcRow=3},
    @{A=8;   Src=$ws1; SrcRow=4},
    @{A=11;  Src=$ws4; SrcRow=2},
    @{A=31;  Src=$ws1; SrcRow=5},
    @{A=37;  Src=$ws3; SrcRow=3},
    @{A=38;  Src=$ws4; SrcRow=3},
    @{A=89;  Src=$ws2; SrcRow=2},
    @{A=158; Src=$ws2; SrcRow=3}
)

$destRow = 2
foreach ($row in $rows) {
    # Copy the whole source row (A:L) so number formats / the bold+border
    # style on column A travel with it, then overwrite column A with the
    # plain index value (copy already wrote the correct one, but be explicit).
    $row.Src.Range("A" + $row.SrcRow + ":L" + $row.SrcRow).Copy($ws5.Range("A" + $destRow + ":L" + $destRow))
    $destRow++
}

Write-Output "all_drugs_present populated"

# ---------------------------------------------------------------
# 2) not_in_DiSCoVER
# ---------------------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws6 = $wb.Worksheets.Add($null, $lastSheet2)
$ws6.Name = "not_in_DiSCoVER"

# Header - reuse the same bold/border/centered style used elsewhere.
$ws1.Range("B1").Copy($ws6.Range("B1"))
$ws6.Range("B1").Value = "not_found"

$ws1.Range("A2").Copy($ws6.Range("A2"))
$ws6.Cells.Item(2, 1).Value = 0
$ws6.Cells.Item(2, 2).Value = "anisomycin"

$ws1.Range("A2").Copy($ws6.Range("A3"))
$ws6.Cells.Item(3, 1).Value = 1
$ws6.Cells.Item(3, 2).Value = "monafide"

Write-Output "not_in_DiSCoVER populated"
